$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for LIDC-IDRI .../tcia-lidc-xml/189/055.xml (row 26),
# which shifts every following row up by one (B189 -> B188).
$ws.Rows(26).Delete()

# Replace the ", " separator used between RID tokens in column B with " | "
# for every data row (row 1 is just the "RIDs" header, no separator there).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace ", ", " | "
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
